# Update the "Week 6 | Lecture 1 (6.1)" header text to "Week 5 | Lecture 3 (5.3)"
# on both the title slide and the duplicated header on the final (notebook) slide.
$p = $ppt.ActivePresentation

$slideIndices = @(1, $p.Slides.Count)

foreach ($idx in $slideIndices) {
    $s = $p.Slides.Item($idx)
    foreach ($sh in $s.Shapes) {
        if ($sh.Name -eq "Subtitle 2" -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "Week 6 | Lecture 1 (6.1)") {
                $tr.Characters(6, 1).Text = "5"
                $tr.Characters(18, 2).Text = "3 "
                $tr.Characters(21, 1).Text = "5"
                $tr.Characters(23, 1).Text = "3"
            }
        }
    }
}
